$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1941.6
$ws.Range("I31").Value = 188
$ws.Range("J31").Value = 6033.3335
$ws.Range("K31").Value = 564
$ws.Range("L31").Value = 18100.0005
$ws.Range("M31").Value = -334
$ws.Range("N31").Value = -18560.0005
$ws.Range("H112").Value = 1682.0952
$ws.Range("I112").Value = 339.66666
$ws.Range("J112").Value = 1905.8334
$ws.Range("K112").Value = 1018.99998
$ws.Range("L112").Value = 5717.5002
$ws.Range("M112").Value = 89.00002000000006
$ws.Range("N112").Value = -7933.5002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1450.4546
$ws.Range("I74").Value = 1389.8
$ws.Range("J74").Value = 2057
$ws.Range("K74").Value = 1389.8
$ws.Range("L74").Value = 2057
$ws.Range("M74").Value = -515.8
$ws.Range("N74").Value = -3805
$ws.Range("H77").Value = 1450.4546
$ws.Range("I77").Value = 1389.8
$ws.Range("J77").Value = 2057
$ws.Range("K77").Value = 6949
$ws.Range("L77").Value = 10285
$ws.Range("M77").Value = -2581
$ws.Range("N77").Value = -19021
$ws.Range("H102").Value = 2362.2
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 2737
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 2737
$ws.Range("M102").Value = -178
$ws.Range("N102").Value = -5981
$ws.Range("H110").Value = 2445.4736
$ws.Range("I110").Value = 1662.5834
$ws.Range("J110").Value = 3787.5715
$ws.Range("K110").Value = 1662.5834
$ws.Range("L110").Value = 3787.5715
$ws.Range("M110").Value = 382.4166
$ws.Range("N110").Value = -7877.5715
$ws.Range("H135").Value = 33705.668
$ws.Range("J135").Value = 33705.668
$ws.Range("L135").Value = 33705.668
$ws.Range("N135").Value = -43845.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1940.7858
$ws.Range("I99").Value = 1250.909
$ws.Range("K99").Value = 1250.909
$ws.Range("M99").Value = 247.0909999999999
$ws.Range("H105").Value = 3069.261
$ws.Range("I105").Value = 2520
$ws.Range("J105").Value = 3151.65
$ws.Range("K105").Value = 2520
$ws.Range("L105").Value = 3151.65
$ws.Range("M105").Value = -773
$ws.Range("N105").Value = -6645.65

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 69175
$ws.Range("J112").Value = 69175
$ws.Range("L112").Value = 69175
$ws.Range("N112").Value = -72129
$ws.Range("H132").Value = 2623.75
$ws.Range("I132").Value = 2401.7778
$ws.Range("J132").Value = 3822.4
$ws.Range("K132").Value = 7205.3334
$ws.Range("L132").Value = 11467.2
$ws.Range("M132").Value = -4675.3334
$ws.Range("N132").Value = -16527.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 506.18518
$ws.Range("I5").Value = 517.2381
$ws.Range("J5").Value = 467.5
$ws.Range("K5").Value = 1551.7143
$ws.Range("L5").Value = 1402.5
$ws.Range("M5").Value = -1439.7143
$ws.Range("N5").Value = -1626.5
$ws.Range("H34").Value = 675.5
$ws.Range("I34").Value = 100
$ws.Range("J34").Value = 739.44446
$ws.Range("K34").Value = 300
$ws.Range("L34").Value = 2218.33338
$ws.Range("M34").Value = -216
$ws.Range("N34").Value = -2386.33338
$ws.Range("H39").Value = 2758.0667
$ws.Range("I39").Value = 600
$ws.Range("J39").Value = 2912.2144
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 8736.643199999999
$ws.Range("M39").Value = -1506
$ws.Range("N39").Value = -9324.643199999999
$ws.Range("H55").Value = 1903.2
$ws.Range("I55").Value = 293.5
$ws.Range("J55").Value = 2976.3333
$ws.Range("K55").Value = 880.5
$ws.Range("L55").Value = 8928.999899999999
$ws.Range("M55").Value = -703.5
$ws.Range("N55").Value = -9282.999899999999
$ws.Range("H109").Value = 3953.6875
$ws.Range("I109").Value = 703
$ws.Range("J109").Value = 4703.846
$ws.Range("K109").Value = 2109
$ws.Range("L109").Value = 14111.538
$ws.Range("M109").Value = -1069
$ws.Range("N109").Value = -16191.538
$ws.Range("H122").Value = 1377987.9
$ws.Range("I122").Value = 466.35715
$ws.Range("J122").Value = 9092108
$ws.Range("K122").Value = 4197.21435
$ws.Range("L122").Value = 81828972
$ws.Range("M122").Value = -1747.21435
$ws.Range("N122").Value = -81833872
$ws.Range("H131").Value = 789.5893
$ws.Range("I131").Value = 446.46155
$ws.Range("J131").Value = 893.32556
$ws.Range("K131").Value = 1339.38465
$ws.Range("L131").Value = 2679.97668
$ws.Range("M131").Value = 3700.61535
$ws.Range("N131").Value = -12759.97668
$ws.Range("H132").Value = 1676.3334
$ws.Range("I132").Value = 322.75
$ws.Range("J132").Value = 2353.125
$ws.Range("K132").Value = 2904.75
$ws.Range("L132").Value = 21178.125
$ws.Range("M132").Value = -374.75
$ws.Range("N132").Value = -26238.125
$ws.Range("H135").Value = 506.18518
$ws.Range("I135").Value = 517.2381
$ws.Range("J135").Value = 467.5
$ws.Range("K135").Value = 4655.142900000001
$ws.Range("L135").Value = 4207.5
$ws.Range("M135").Value = -2120.142900000001
$ws.Range("N135").Value = -9277.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 14466.2
$ws.Range("I35").Value = 777
$ws.Range("K35").Value = 777
$ws.Range("M35").Value = -441
$ws.Range("H57").Value = 2520.5
$ws.Range("I57").Value = 2520.5
$ws.Range("K57").Value = 2520.5
$ws.Range("M57").Value = -1954.5
$ws.Range("H122").Value = 4914.5864
$ws.Range("I122").Value = 4674.913
$ws.Range("J122").Value = 5833.3335
$ws.Range("K122").Value = 14024.739
$ws.Range("L122").Value = 17500.0005
$ws.Range("M122").Value = -11574.739
$ws.Range("N122").Value = -22400.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 24000
$ws.Range("J109").Value = 24000
$ws.Range("L109").Value = 24000
$ws.Range("N109").Value = -26774
$ws.Range("H122").Value = 1528.4857
$ws.Range("I122").Value = 1419.96
$ws.Range("J122").Value = 1799.8
$ws.Range("K122").Value = 4259.88
$ws.Range("L122").Value = 5399.4
$ws.Range("M122").Value = -1809.88
$ws.Range("N122").Value = -10299.4
$ws.Range("H132").Value = 2056.3
$ws.Range("I132").Value = 1014.6
$ws.Range("J132").Value = 3792.4666
$ws.Range("K132").Value = 3043.8
$ws.Range("L132").Value = 11377.3998
$ws.Range("M132").Value = -513.8000000000002
$ws.Range("N132").Value = -16437.3998
